$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 and G3 values
$ws.Range("G2").Value = 60
$ws.Range("G3").Value = 60

# Update G4:G9 values and apply number format (#,##0)
$ws.Range("G4:G9").Value = 52.961
$ws.Range("G4:G9").NumberFormat = "#,##0"

# Update the selection to G4:G9 with active cell G4
$ws.Range("G4:G9").Select()
